# Added OK endpoints for DAC Wishbones
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DAC80508")

# New register rows for the DAC Wishbone OK endpoints.
# Row 19 (WB_IN_1) is entered before row 18 (WB_IN_0) to reproduce the
# original authoring order of the new shared strings.
$ws.Cells.Item(19, 1).Value = "WB_IN_1"
$ws.Cells.Item(19, 2).Value = "0x04"
$ws.Cells.Item(19, 3).Value = "0x00000000"
$ws.Cells.Item(19, 4).Value = 32
$ws.Cells.Item(19, 5).Value = "None"
$ws.Cells.Item(19, 6).Value = "None"

$ws.Cells.Item(18, 1).Value = "WB_IN_0"
$ws.Cells.Item(18, 2).Value = "0x03"
$ws.Cells.Item(18, 3).Value = "0x00000000"
$ws.Cells.Item(18, 4).Value = 32
$ws.Cells.Item(18, 5).Value = "None"
$ws.Cells.Item(18, 6).Value = "None"

$ws.Cells.Item(20, 1).Value = "WB_OUT_0"
$ws.Cells.Item(20, 2).Value = "0x22"
$ws.Cells.Item(20, 3).Value = "0x00000000"
$ws.Cells.Item(20, 4).Value = 32
$ws.Cells.Item(20, 5).Value = "None"
$ws.Cells.Item(20, 6).Value = "None"

$ws.Cells.Item(21, 1).Value = "WB_OUT_1"
$ws.Cells.Item(21, 2).Value = "0x23"
$ws.Cells.Item(21, 3).Value = "0x00000000"
$ws.Cells.Item(21, 4).Value = 32
$ws.Cells.Item(21, 5).Value = "None"
$ws.Cells.Item(21, 6).Value = "None"

# Move the active tab / window selection onto the DAC80508 sheet
$ws.Activate()
$ws.Range("A22").Select() | Out-Null
